$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.02406616871877758
$ws.Range("J2").Value = 0.02406616871877757
$ws.Range("M2").Value = 15.42521333333333
$ws.Range("N2").Value = 46.27564
$ws.Range("O2").Value = 0.2413167708794786
$ws.Range("P2").Value = 0.2413167708794786
$ws.Range("Q2").Value = 0.6913374946488889
$ws.Range("R2").Value = 6.22203745184
$ws.Range("S2").Value = 0.005807570122656125
$ws.Range("T2").Value = 0.005807570122656123
$ws.Range("I3").Value = 0.02406616871877758
$ws.Range("J3").Value = 0.02406616871877757
$ws.Range("O3").Value = 0.2555517499330554
$ws.Range("P3").Value = 0.2555517499330554
$ws.Range("S3").Value = 0.006150151530267768
$ws.Range("T3").Value = 0.006150151530267767
$ws.Range("I4").Value = 0.02406616871877758
$ws.Range("J4").Value = 0.02406616871877757
$ws.Range("M4").Value = 15.93058
$ws.Range("N4").Value = 47.79174
$ws.Range("O4").Value = 0.2492228820932917
$ws.Range("P4").Value = 0.2492228820932917
$ws.Range("Q4").Value = 0.7139873548266668
$ws.Range("R4").Value = 6.42588619344
$ws.Range("S4").Value = 0.00599783992903717
$ws.Range("T4").Value = 0.005997839929037168
$ws.Range("I5").Value = 0.02406616871877758
$ws.Range("J5").Value = 0.02406616871877757
$ws.Range("M5").Value = 4.108632
$ws.Range("N5").Value = 12.325896
$ws.Range("O5").Value = 0.06427669981260728
$ws.Range("P5").Value = 0.06427669981260727
$ws.Range("Q5").Value = 0.184143408064
$ws.Range("R5").Value = 1.657290672576
$ws.Range("S5").Value = 0.001546893902376426
$ws.Range("T5").Value = 0.001546893902376425
$ws.Range("I6").Value = 0.02406616871877758
$ws.Range("J6").Value = 0.02406616871877757
$ws.Range("M6").Value = 12.12146366666667
$ws.Range("N6").Value = 36.364391
$ws.Range("O6").Value = 0.189631897281567
$ws.Range("P6").Value = 0.189631897281567
$ws.Range("Q6").Value = 0.5432678395884444
$ws.Range("R6").Value = 4.889410556295999
$ws.Range("S6").Value = 0.00456371323444009
$ws.Range("T6").Value = 0.004563713234440089
$ws.Range("I7").Value = 0.8626970447097064
$ws.Range("J7").Value = 0.8626970447097063
$ws.Range("M7").Value = 15.42521333333333
$ws.Range("N7").Value = 46.27564
$ws.Range("O7").Value = 0.2413167708794786
$ws.Range("P7").Value = 0.2413167708794786
$ws.Range("Q7").Value = 24.78229170999111
$ws.Range("R7").Value = 223.04062538992
$ws.Range("S7").Value = 0.2081832650766156
$ws.Range("T7").Value = 0.2081832650766155
$ws.Range("I8").Value = 0.8626970447097064
$ws.Range("J8").Value = 0.8626970447097063
$ws.Range("O8").Value = 0.2555517499330554
$ws.Range("P8").Value = 0.2555517499330554
$ws.Range("S8").Value = 0.2204637394376408
$ws.Range("T8").Value = 0.2204637394376408
$ws.Range("I9").Value = 0.8626970447097064
$ws.Range("J9").Value = 0.8626970447097063
$ws.Range("M9").Value = 15.93058
$ws.Range("N9").Value = 47.79174
$ws.Range("O9").Value = 0.2492228820932917
$ws.Range("P9").Value = 0.2492228820932917
$ws.Range("Q9").Value = 25.59421851341333
$ws.Range("R9").Value = 230.34796662072
$ws.Range("S9").Value = 0.2150038438559184
$ws.Range("T9").Value = 0.2150038438559183
$ws.Range("I10").Value = 0.8626970447097064
$ws.Range("J10").Value = 0.8626970447097063
$ws.Range("M10").Value = 4.108632
$ws.Range("N10").Value = 12.325896
$ws.Range("O10").Value = 0.06427669981260728
$ws.Range("P10").Value = 0.06427669981260727
$ws.Range("Q10").Value = 6.600966518431999
$ws.Range("R10").Value = 59.40869866588799
$ws.Range("S10").Value = 0.05545131897202924
$ws.Range("T10").Value = 0.05545131897202922
$ws.Range("I11").Value = 0.8626970447097064
$ws.Range("J11").Value = 0.8626970447097063
$ws.Range("M11").Value = 12.12146366666667
$ws.Range("N11").Value = 36.364391
$ws.Range("O11").Value = 0.189631897281567
$ws.Range("P11").Value = 0.189631897281567
$ws.Range("Q11").Value = 19.47445666052755
$ws.Range("R11").Value = 175.270109944748
$ws.Range("S11").Value = 0.1635948773675024
$ws.Range("T11").Value = 0.1635948773675024
$ws.Range("E12").Value = 3.0
$ws.Range("F12").Value = 1.0
$ws.Range("G12").Value = 0.210882
$ws.Range("H12").Value = 0.6326459999999999
$ws.Range("I12").Value = 0.113236786571516
$ws.Range("J12").Value = 0.113236786571516
$ws.Range("M12").Value = 15.42521333333333
$ws.Range("N12").Value = 46.27564
$ws.Range("O12").Value = 0.2413167708794786
$ws.Range("P12").Value = 0.2413167708794786
$ws.Range("Q12").Value = 3.25289983816
$ws.Range("R12").Value = 29.27609854344
$ws.Range("S12").Value = 0.02732593568020696
$ws.Range("T12").Value = 0.02732593568020695
$ws.Range("E13").Value = 3.0
$ws.Range("F13").Value = 1.0
$ws.Range("G13").Value = 0.210882
$ws.Range("H13").Value = 0.6326459999999999
$ws.Range("I13").Value = 0.113236786571516
$ws.Range("J13").Value = 0.113236786571516
$ws.Range("O13").Value = 0.2555517499330554
$ws.Range("P13").Value = 0.2555517499330554
$ws.Range("Q13").Value = 3.444784392602
$ws.Range("R13").Value = 31.003059533418
$ws.Range("S13").Value = 0.02893785896514683
$ws.Range("T13").Value = 0.02893785896514682
$ws.Range("E14").Value = 3.0
$ws.Range("F14").Value = 1.0
$ws.Range("G14").Value = 0.210882
$ws.Range("H14").Value = 0.6326459999999999
$ws.Range("I14").Value = 0.113236786571516
$ws.Range("J14").Value = 0.113236786571516
$ws.Range("M14").Value = 15.93058
$ws.Range("N14").Value = 47.79174
$ws.Range("O14").Value = 0.2492228820932917
$ws.Range("P14").Value = 0.2492228820932917
$ws.Range("Q14").Value = 3.35947257156
$ws.Range("R14").Value = 30.23525314404
$ws.Range("S14").Value = 0.02822119830833618
$ws.Range("T14").Value = 0.02822119830833617
$ws.Range("E15").Value = 3.0
$ws.Range("F15").Value = 1.0
$ws.Range("G15").Value = 0.210882
$ws.Range("H15").Value = 0.6326459999999999
$ws.Range("I15").Value = 0.113236786571516
$ws.Range("J15").Value = 0.113236786571516
$ws.Range("M15").Value = 4.108632
$ws.Range("N15").Value = 12.325896
$ws.Range("O15").Value = 0.06427669981260728
$ws.Range("P15").Value = 0.06427669981260727
$ws.Range("Q15").Value = 0.866436533424
$ws.Range("R15").Value = 7.797928800815999
$ws.Range("S15").Value = 0.007278486938201615
$ws.Range("T15").Value = 0.007278486938201612
$ws.Range("E16").Value = 3.0
$ws.Range("F16").Value = 1.0
$ws.Range("G16").Value = 0.210882
$ws.Range("H16").Value = 0.6326459999999999
$ws.Range("I16").Value = 0.113236786571516
$ws.Range("J16").Value = 0.113236786571516
$ws.Range("M16").Value = 12.12146366666667
$ws.Range("N16").Value = 36.364391
$ws.Range("O16").Value = 0.189631897281567
$ws.Range("P16").Value = 0.189631897281567
$ws.Range("Q16").Value = 2.556198500953999
$ws.Range("R16").Value = 23.005786508586
$ws.Range("S16").Value = 0.02147330667962445
$ws.Range("T16").Value = 0.02147330667962445
